$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column R (the last existing data column) into column S
$ws.Range("R2:R6").Copy()
$ws.Range("S2:S6").PasteSpecial(-4122)  # xlPasteFormats

# Set new values for the new column S (year 2022 data)
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 265803
$ws.Range("S5").Value = 3.8
$ws.Range("S6").Value = 33.6

# Update the selected cell to match the new view state
$ws.Range("C19").Select()
